$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet had a leading "row index" column (A) that is no longer wanted.
# Deleting it shifts every other column (B..F) one place to the left
# (B->A, C->B, D->C, E->D, F->E), which is exactly the transformation
# described by the diff: header row loses its old A-column formatting,
# and the remaining data keeps its values/styles, now one column over.
$ws.Columns.Item(1).Delete()
